# DISCOVERYACCESS-7433: Add hotel (Nestle Library) permanent reserve to the
# hierarchical location facet mapping sheet.
#
# A new row is inserted above the existing "Nestle Library Reserve" row
# (row 72) so the new "Permanent Reserve" mapping sits right next to the
# related Nestle Library entries, and every following row shifts down by
# one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 72, pushing rows 72-118 down to 73-119.
$ws.Rows(72).Insert()

# Populate the new row - Voyager Display Name (A) and Facet Display Name (E).
# Set E before A so the shared-string table picks up the same ordering as
# the authored workbook (the " > " facet name first, then the plain name).
$ws.Range("E72").Value = "Nestle Library > Permanent Reserve"
$ws.Range("A72").Value = "Nestle Library Permanent Reserve"

# Match the author's on-screen selection after making the edit.
[void]$ws.Range("A72").Select()
